$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F ("dSF") values to re-pulled / re-pushed data per repull.
$ws.Range("F2").Value = -2
$ws.Range("F3").Value = -9
$ws.Range("F6").Value = -10
$ws.Range("F7").Value = -1
$ws.Range("F11").Value = -4
$ws.Range("F13").Value = -5
$ws.Range("F14").Value = 2
$ws.Range("F17").Value = -7
$ws.Range("F18").Value = -6
$ws.Range("F27").Value = -5
